# "atualizacao da funcao de comprimentos"
#
# Adds three new columns to Sheet1 that normalise the existing
# volume / area / comprimento (length) measurements to a standard
# 1 dm3 soil-core volume:
#   I: volume_cm3_dm3      = volume (E)      / 475
#   J: area_cm2_dm3        = area (F)        / 47.5
#   K: comprimento_cm_dm3  = comprimento (H) / 4.75

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1), matching the style of the existing headers.
$ws.Range("I1").Value = "volume_cm3_dm3"
$ws.Range("J1").Value = "area_cm2_dm3"
$ws.Range("K1").Value = "comprimento_cm_dm3"

# Data rows 2..110: derive the normalised columns from the existing
# volume (E), area (F) and comprimento (H) columns.
$lastRow = 110
for ($r = 2; $r -le $lastRow; $r++) {
    $volume = $ws.Range("E$r").Value()
    $area = $ws.Range("F$r").Value()
    $comprimento = $ws.Range("H$r").Value()

    $ws.Range("I$r").Value = $volume / 475
    $ws.Range("J$r").Value = $area / 47.5
    $ws.Range("K$r").Value = $comprimento / 4.75
}
